# Apply cryptocurrency price/volume updates from the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.324.34'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.986.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.92%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.07'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.47%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.986.19'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.81%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.69%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.07'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.435'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000222'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.47'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -7.31%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.473.05'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.271.71'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.979.95'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.14'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.33'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -7.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.03'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.658'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.12'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -6.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.90'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.81'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.50'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.10'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -8.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.88'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -7.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.13'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -7.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.00'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -11.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0923'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -10.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.25'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.943'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -9.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.52'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '49.39'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0641'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -9.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0356'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.82'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.107'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '372.95'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.641.42'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.37'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -8.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.233'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.76'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.95'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -8.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.106'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.45%  '
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '32.07'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.15%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.27'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -8.41%  '
